$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "obs" column value text (shared string used by K2:K54)
$ws.Range("K2:K54").Value = "dN/N_inel-dy-dpT"

# Widen column K
$ws.Columns.Item(11).ColumnWidth = 16.25

# Change selection to K2
$ws.Range("K2").Select()
